$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the team member field (B5 merged label is static; the blank area after
# the colon now carries the team member's name)
$ws.Range("B5").Value = "MEMBRES DE L'EQUIPE : AMMAR-BOUDJELAL Lina"

# Update the date field
$ws.Range("E4").Value = "date : 23/05/2025"

# Fill in the project name cell (merged C4:D4), previously blank
$ws.Range("C4").Value = "ENSemenC"

# Update the selection to match the authored state
$ws.Range("C4:D4").Select()
